$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 66: Running Sum of 1d Array ---
$ws.Range("B66").Value = "Prefix Sum"
$ws.Range("C66").Value = "Running Sum of 1d Array"
$ws.Range("D66").Value = "Easy"
$ws.Range("E66").Value = "Done"
$ws.Range("F65").Copy()
$ws.Range("F66").PasteSpecial(-4122)
$ws.Range("F66").Value = 45910
$ws.Range("G66").Value = "O(n)"
$ws.Range("H66").Value = "O(n)"
$ws.Range("I66").Value = "Prefix Sum"

# --- Row 67: Find Pivot Index ---
$ws.Range("B67").Value = "Prefix Sum"
$ws.Range("C67").Value = "Find Pivot Index"
$ws.Range("D67").Value = "Easy"
$ws.Range("E67").Value = "Done"
$ws.Range("F65").Copy()
$ws.Range("F67").PasteSpecial(-4122)
$ws.Range("F67").Value = 45910
$ws.Range("G67").Value = "O(n)"
$ws.Range("H67").Value = "O(1)"
$ws.Range("I67").Value = "Prefix Sum"

# --- Row 68: Subarray Sum Equals K ---
$ws.Range("B68").Value = "Prefix Sum"
$ws.Range("C68").Value = "Subarray Sum Equals K"
$ws.Range("D68").Value = "Medium"
$ws.Range("E68").Value = "Done"
$ws.Range("F65").Copy()
$ws.Range("F68").PasteSpecial(-4122)
$ws.Range("F68").Value = 45910
$ws.Range("G68").Value = "O(n)"
$ws.Range("H68").Value = "O(n)"

# --- Row 69: Range Sum Query - Immutable (written before "Prefix Sum + HashMap" below,
#             to match shared-string append order) ---
$ws.Range("B69").Value = "Prefix Sum"
$ws.Range("C69").Value = "Range Sum Query " + [char]8211 + " Immutable"
$ws.Range("D69").Value = "Easy"
$ws.Range("E69").Value = "Done"
$ws.Range("G69").Value = "O(n)"
$ws.Range("H69").Value = "O(n)"
$ws.Range("I69").Value = "Prefix Sum"

# "Prefix Sum + HashMap" is introduced here (after "Range Sum Query - Immutable")
$ws.Range("I68").Value = "Prefix Sum + HashMap"

$excel.CutCopyMode = 0

# --- Update sheet view scroll/selection ---
$ws.Application.ActiveWindow.ScrollRow = 56
$ws.Range("G70").Select()
